# "Generate Report for Handback"
#
# The localization-status report is refreshed after a handback run:
#   - the status text moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + per-language sheets)
#   - the per-language sheets (zh-cn, de-de) gain the generated target
#     (.md) and handback (.xlf) file links for each row, via real
#     hyperlinks in the "Latest Target File" (I) / "Latest Handback
#     File" (J) columns
#   - de-de (which has actually been handed back) also stamps the
#     "Latest Handback DateTime" (K) column with the real timestamp;
#     zh-cn's placeholder timestamp text is refreshed as well
#   - a few columns are widened so the new/longer text fits

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared across Overview!E:F and the zh-cn/de-de Status column)
# ---------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# ---------------------------------------------------------------------
# 2. Column widths
#    ColumnWidth is in characters; the saved OOXML width is
#    ColumnWidth + 5/6. Target widths of 29.9777047293527 and 40 are
#    reached (as closely as the engine's internal pixel snapping
#    allows) with these inputs.
# ---------------------------------------------------------------------
$wStatus = 175 / 6   # -> saved width 30   (closest reachable to 29.9777047293527)
$wWide40 = 235 / 6   # -> saved width 40   (exact)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $wStatus   # E: zh-cn
$wsOverview.Columns.Item(6).ColumnWidth = $wStatus   # F: de-de

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = $wStatus         # C: Status
$wsZh.Columns.Item(9).ColumnWidth = $wWide40         # I: Latest Target File
$wsZh.Columns.Item(10).ColumnWidth = $wWide40        # J: Latest Handback File

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = $wStatus         # C: Status
$wsDe.Columns.Item(9).ColumnWidth = $wWide40         # I: Latest Target File
$wsDe.Columns.Item(10).ColumnWidth = $wWide40        # J: Latest Handback File

# ---------------------------------------------------------------------
# 3. Populate "Latest Target File" (I) and "Latest Handback File" (J)
#    with real hyperlinks for both tracked source files, on both
#    language sheets.
# ---------------------------------------------------------------------
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/621ae638fdf405cea936e34bc356a9f494abfea4/e2e/9f148c15-ae5b-4119-8944-6d7ea72aaf17.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/621ae638fdf405cea936e34bc356a9f494abfea4/e2e/cb4dc240-bbcf-4a5b-8475-8f268fea9a70.md"
$mdDisp1 = "9f148c15-ae5b-4119-8944-6d7ea72aaf17.md"
$mdDisp2 = "cb4dc240-bbcf-4a5b-8475-8f268fea9a70.md"

# zh-cn
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdDisp1)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdDisp2)
$wsZh.Range("J2").Value = "9f148c15-ae5b-4119-8944-6d7ea72aaf17.2def155db845fc158ff9ecda8448134afb27b4bb.zh-cn.xlf"
$wsZh.Range("J3").Value = "cb4dc240-bbcf-4a5b-8475-8f268fea9a70.a304c472ecf48fedafb0bc83a52c245b85010114.zh-cn.xlf"

# de-de
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdDisp1)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdDisp2)
$wsDe.Range("J2").Value = "9f148c15-ae5b-4119-8944-6d7ea72aaf17.2def155db845fc158ff9ecda8448134afb27b4bb.de-de.xlf"
$wsDe.Range("J3").Value = "cb4dc240-bbcf-4a5b-8475-8f268fea9a70.a304c472ecf48fedafb0bc83a52c245b85010114.de-de.xlf"

# ---------------------------------------------------------------------
# 4. Latest Handback DateTime (K)
#    de-de has actually been handed back -> real timestamp.
#    zh-cn hasn't, but its placeholder text is refreshed too.
# ---------------------------------------------------------------------
$wsDe.Range("K2").Value = "2016-08-30 06:50:46"
$wsDe.Range("K3").Value = "2016-08-30 06:50:46"

$wsZh.Cells.Replace("0001-01-01 00:00:00", "2016-08-30 06:50:39")
